$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with default (unstyled) formatting, used to reset style
# after forcing numeric-looking text into Price cells so they stay as text
# without leaving a stray style index on the cell.
$cleanStyle = $ws.Range("E2").Style

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.446.69'
$ws.Range('D2').Style = $cleanStyle
$ws.Range('E2').Value = '  +0.40%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.871.38'
$ws.Range('D3').Style = $cleanStyle
$ws.Range('E3').Value = '  -0.40%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('D4').Style = $cleanStyle
$ws.Range('E4').Value = '  -0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7074'
$ws.Range('D5').Style = $cleanStyle
$ws.Range('E5').Value = '  -0.39%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '243.66'
$ws.Range('D6').Style = $cleanStyle
$ws.Range('E6').Value = '  +0.51%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9998'
$ws.Range('D7').Style = $cleanStyle
$ws.Range('E7').Value = '  -0.09%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3160'
$ws.Range('D8').Style = $cleanStyle
$ws.Range('E8').Value = '  +0.64%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07863'
$ws.Range('D9').Style = $cleanStyle
$ws.Range('E9').Value = '  -1.95%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.51'
$ws.Range('D10').Style = $cleanStyle
$ws.Range('E10').Value = '  -2.43%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07997'
$ws.Range('D11').Style = $cleanStyle
$ws.Range('E11').Value = '  -4.03%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.863.43'
$ws.Range('D12').Style = $cleanStyle
$ws.Range('E12').Value = '  -1.14%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.225'
$ws.Range('D13').Style = $cleanStyle
$ws.Range('E13').Value = '  -0.96%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '94.20'
$ws.Range('D14').Style = $cleanStyle
$ws.Range('E14').Value = '  -0.67%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.7039'
$ws.Range('D15').Style = $cleanStyle
$ws.Range('E15').Value = '  -1.90%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.495'
$ws.Range('D16').Style = $cleanStyle
$ws.Range('E16').Value = '  +1.84%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.451.70'
$ws.Range('D17').Style = $cleanStyle
$ws.Range('E17').Value = '  +0.35%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008350'
$ws.Range('D18').Style = $cleanStyle
$ws.Range('E18').Value = '  -4.44%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '257.23'
$ws.Range('D19').Style = $cleanStyle
$ws.Range('E19').Value = '  +5.94%  '

# Row 20
$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.120.81'
$ws.Range('D20').Style = $cleanStyle
$ws.Range('E20').Value = '  -0.64%  '

# Row 21
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.21'
$ws.Range('D21').Style = $cleanStyle
$ws.Range('E21').Value = '  -0.63%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9999'
$ws.Range('D22').Style = $cleanStyle
$ws.Range('E22').Value = '  +0.01%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.631'
$ws.Range('D23').Style = $cleanStyle
$ws.Range('E23').Value = '  -2.78%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9995'
$ws.Range('D24').Style = $cleanStyle
$ws.Range('E24').Value = '  -0.12%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1561'
$ws.Range('D25').Style = $cleanStyle
$ws.Range('E25').Value = '  -0.74%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.072'
$ws.Range('D26').Style = $cleanStyle
$ws.Range('E26').Value = '  +0.05%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '160.62'
$ws.Range('D27').Style = $cleanStyle
$ws.Range('E27').Value = '  -1.79%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.92'
$ws.Range('D28').Style = $cleanStyle
$ws.Range('E28').Value = '  +1.71%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.501'
$ws.Range('D29').Style = $cleanStyle
$ws.Range('E29').Value = '  -0.40%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.340'
$ws.Range('D30').Style = $cleanStyle
$ws.Range('E30').Value = '  -1.99%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.259'
$ws.Range('D31').Style = $cleanStyle
$ws.Range('E31').Value = '  -1.99%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.206'
$ws.Range('D32').Style = $cleanStyle
$ws.Range('E32').Value = '  +0.03%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05315'
$ws.Range('D33').Style = $cleanStyle
$ws.Range('E33').Value = '  -1.88%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.901'
$ws.Range('D34').Style = $cleanStyle
$ws.Range('E34').Value = '  -1.98%  '

# Row 35
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.175'
$ws.Range('D35').Style = $cleanStyle
$ws.Range('E35').Value = '  -0.34%  '

# Row 36
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.712'
$ws.Range('D36').Style = $cleanStyle
$ws.Range('E36').Value = '  -3.63%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.712'
$ws.Range('D37').Style = $cleanStyle
$ws.Range('E37').Value = '  +0.89%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01884'
$ws.Range('D38').Style = $cleanStyle
$ws.Range('E38').Value = '  -0.10%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.267.51'
$ws.Range('D39').Style = $cleanStyle
$ws.Range('E39').Value = '  -0.66%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.753'
$ws.Range('D40').Style = $cleanStyle
$ws.Range('E40').Value = '  +0.33%  '

# Row 41
$ws.Range('E41').Value = '  -1.58%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '108.88'
$ws.Range('D42').Style = $cleanStyle
$ws.Range('E42').Value = '  -3.46%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.996'
$ws.Range('D43').Style = $cleanStyle
$ws.Range('E43').Value = '  -8.50%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '71.61'
$ws.Range('D44').Style = $cleanStyle
$ws.Range('E44').Value = '  -3.98%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9994'
$ws.Range('D45').Style = $cleanStyle
$ws.Range('E45').Value = '  -0.08%  '

# Row 46
$ws.Range('E46').Value = '  +2.25%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.019.84'
$ws.Range('D47').Style = $cleanStyle
$ws.Range('E47').Value = '  -0.54%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5192'
$ws.Range('D48').Style = $cleanStyle
$ws.Range('E48').Value = '  -0.48%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.792'
$ws.Range('D49').Style = $cleanStyle
$ws.Range('E49').Value = '  -0.81%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.533'
$ws.Range('D50').Style = $cleanStyle
$ws.Range('E50').Value = '  -0.09%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4331'
$ws.Range('D51').Style = $cleanStyle
$ws.Range('E51').Value = '  -1.02%  '
